$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vehicles")

$ws.Range("B2").Value = "N,m,kg,s,C"
$ws.Range("C2").Value = "0|2|5|2"
$ws.Range("B3").Value = "N,m,kg,s,C"
$ws.Range("C3").Value = "0|2|5|2|0"
